# Update database (income statement) + shift each metric one period to the
# left and append the newly published period (1401/12) in column H.
# Mirrors: drop oldest "1396/12" column, shift 1397..1400 left one column,
# add new "1401/12" column with freshly reported figures; publish-date row
# and figures shift/refresh the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: "دوره مالی" (fiscal period) headers ---------------------------
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" (publish date) ----------------------------------
$ws.Range("D9").Value = "1399-01-24 (7)"
$ws.Range("E9").Value = "1400-02-04 (7)"
$ws.Range("F9").Value = "1401-01-31 (8)"
$ws.Range("G9").Value = "1402-01-30 (9)"
$ws.Range("H9").Value = "1402-01-30 (2)"

# --- Row 11: فروش (sales) ---------------------------------------------------
$ws.Range("D11").Value = 16814
$ws.Range("E11").Value = 17551
$ws.Range("F11").Value = 19221
$ws.Range("G11").Value = 28449
$ws.Range("H11").Value = 28949

# --- Row 12: بهای تمام شده کالای فروش رفته (COGS) ---------------------------
$ws.Range("D12").Value = -11056
$ws.Range("E12").Value = -10214
$ws.Range("F12").Value = -9341
$ws.Range("G12").Value = -12338
$ws.Range("H12").Value = -13109

# --- Row 13: سود (زیان) ناخالص (gross profit) -------------------------------
$ws.Range("D13").Value = 5758
$ws.Range("E13").Value = 7336
$ws.Range("F13").Value = 9881
$ws.Range("G13").Value = 16111
$ws.Range("H13").Value = 15841

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses) --------------
$ws.Range("D14").Value = -747
$ws.Range("E14").Value = -492
$ws.Range("F14").Value = -377
$ws.Range("G14").Value = -1080
$ws.Range("H14").Value = -1404

# Row 15 (هزینه کاهش ارزش دریافتنی‌ها) is unchanged, all "-"

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی -------------------------
$ws.Range("D16").Value = 198
$ws.Range("E16").Value = 327
$ws.Range("F16").Value = 222
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = 1044

# --- Row 17: سود (زیان) عملیاتی (operating income) --------------------------
$ws.Range("D17").Value = 5209
$ws.Range("E17").Value = 7171
$ws.Range("F17").Value = 9726
$ws.Range("G17").Value = 15045
$ws.Range("H17").Value = 15480

# --- Row 18: هزینه های مالی (finance costs) - mix of numbers and "-" -------
$ws.Range("D18").Value = -23
$ws.Range("E18").Value = "-"
$ws.Range("F18").Value = -24
$ws.Range("G18").Value = "-"
$ws.Range("H18").Value = -85

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی -----------------------
$ws.Range("D19").Value = 200
$ws.Range("E19").Value = 535
$ws.Range("F19").Value = 2621
$ws.Range("G19").Value = 3675
$ws.Range("H19").Value = 1697

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات -------------
$ws.Range("D20").Value = 5385
$ws.Range("E20").Value = 7706
$ws.Range("F20").Value = 12323
$ws.Range("G20").Value = 18720
$ws.Range("H20").Value = 17092

# --- Row 21: مالیات (tax) ----------------------------------------------------
$ws.Range("D21").Value = -831
$ws.Range("E21").Value = -560
$ws.Range("F21").Value = -900
$ws.Range("G21").Value = -1667
$ws.Range("H21").Value = -1383

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ----------------------------
$ws.Range("D22").Value = 4554
$ws.Range("E22").Value = 7146
$ws.Range("F22").Value = 11423
$ws.Range("G22").Value = 17052
$ws.Range("H22").Value = 15709

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی - now has values
$ws.Range("D23").Value = "-"
$ws.Range("E23").Value = "-"
$ws.Range("F23").Value = "-"
$ws.Range("G23").Value = 24
$ws.Range("H23").Value = 16

# --- Row 24: سود (زیان) خالص (net income) -----------------------------------
$ws.Range("D24").Value = 4554
$ws.Range("E24").Value = 7146
$ws.Range("F24").Value = 11423
$ws.Range("G24").Value = 17052
$ws.Range("H24").Value = 15709

# Row 25 (سود هر سهم پس از کسر مالیات) unchanged, all zeros

# --- Row 26: سرمایه (capital) ------------------------------------------------
$ws.Range("D26").Value = 6920
$ws.Range("E26").Value = 5457
$ws.Range("F26").Value = 3096
$ws.Range("G26").Value = 2653
$ws.Range("H26").Value = 2834

# Row 27 (سود هر سهم بر اساس آخرین سرمایه) unchanged, all zeros
